$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.733.57'
$ws.Range('E2').Value = '  -2.99%  '
$ws.Range('D3').Value = '1.785.35'
$ws.Range('E3').Value = '  -2.94%  '
$ws.Range('E4').Value = '  +0.45%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '241.62'
$ws.Range('E5').Value = '  -7.36%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.005'
$ws.Range('E6').Value = '  +0.42%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5078'
$ws.Range('E7').Value = '  -3.40%  '
$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '42.44'
$ws.Range('E8').Value = '  -5.00%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2507'
$ws.Range('E9').Value = '  -21.43%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06150'
$ws.Range('E10').Value = '  -9.50%  '
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').Value = '1.815.88'
$ws.Range('E11').Value = '  -1.27%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.06892'
$ws.Range('E12').Value = '  -11.15%  '
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '14.97'
$ws.Range('E13').Value = '  -20.30%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.6166'
$ws.Range('E14').Value = '  -21.40%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '78.60'
$ws.Range('E15').Value = '  -10.62%  '
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.414'
$ws.Range('E16').Value = '  -11.99%  '
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.004'
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.005'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '25.777.46'
$ws.Range('E19').Value = '  -2.88%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.54'
$ws.Range('E20').Value = '  -16.75%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.044.54'
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('B22').Value = 'ShibaInu'
$ws.Range('C22').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.000006259'
$ws.Range('E22').Value = '  -21.22%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.950'
$ws.Range('E23').Value = '  -14.68%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.231'
$ws.Range('E24').Value = '  -12.46%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '8.060'
$ws.Range('E25').Value = '  -13.72%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '131.70'
$ws.Range('E26').Value = '  -7.22%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.892'
$ws.Range('E27').Value = '  -13.91%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '14.56'
$ws.Range('E28').Value = '  -14.14%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.358'
$ws.Range('E29').Value = '  -19.06%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '98.27'
$ws.Range('E30').Value = '  -12.15%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08315'
$ws.Range('E31').Value = '  -4.47%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.611'
$ws.Range('E32').Value = '  -13.39%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.04339'
$ws.Range('E33').Value = '  -11.04%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.742'
$ws.Range('E34').Value = '  -4.12%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.152'
$ws.Range('E35').Value = '  -22.77%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.038'
$ws.Range('E36').Value = '  -8.66%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.6283'
$ws.Range('E37').Value = '  -13.92%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.848'
$ws.Range('E38').Value = '  -7.93%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.087'
$ws.Range('E39').Value = '  -7.07%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.006'
$ws.Range('E40').Value = '  +0.50%  '
$ws.Range('B41').Value = 'PaxosStandard'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.005'
$ws.Range('E41').Value = '  +0.42%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '101.20'
$ws.Range('E42').Value = '  -7.70%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.01459'
$ws.Range('E43').Value = '  -16.89%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.3911'
$ws.Range('E44').Value = '  -18.79%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.7629'
$ws.Range('E45').Value = '  -14.73%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '5.182'
$ws.Range('E46').Value = '  -12.72%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '6.211'
$ws.Range('E47').Value = '  -18.94%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.05264'
$ws.Range('E48').Value = '  -10.07%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '52.91'
$ws.Range('E49').Value = '  -11.34%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.006'
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '29.39'
$ws.Range('E51').Value = '  -15.86%  '
